$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.908.16"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.462.95"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "2.464.15"
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.111"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "2.910.57"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "62.799.29"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "2.460.45"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "327.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("E22").Value = "  +10.64%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +21.21%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "659.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.02%  "
$ws.Range("D28").Value = "0.0₃0982"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -14.50%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.134"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.66%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +4.05%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "151.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("E42").Value = "  +1.73%  "
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("E44").Value = "  -31.16%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "153.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.23"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.65"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.607"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "  -0.39%  "

Write-Host "Updated cryptos list"
